$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add the new "Keys" (people) first, so the new shared strings are
# inserted ahead of the (about to be edited) message string, matching
# the target shared-string ordering. ---
$ws.Range("B4").Value = "Isis"
$ws.Range("B5").Value = "Viviane"
$ws.Range("B6").Value = "Renata"
$ws.Range("B7").Value = "Alessandra"

# --- Update the message text (shared by every row in column D) ---
$msg = "Nós do @merendinhajf gostariámos de conhecer melhor você e seu pequeno(a). Preencha o formulario pelo link abaixo para nos ajudar a cada vez mais atende-los melhor.  **Mensagem teste do BOT Automaizador de envio de mensagens no Whatsapp**"

$ws.Range("D2").Value = $msg
$ws.Range("D3").Value = $msg
$ws.Range("D4").Value = $msg
$ws.Range("D5").Value = $msg
$ws.Range("D6").Value = $msg
$ws.Range("D7").Value = $msg

# --- Fill in the rest of the new rows (number + phone) ---
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = 5532988329968

$ws.Range("A5").Value = 4
$ws.Range("C5").Value = 5532988090045

$ws.Range("A6").Value = 5
$ws.Range("C6").Value = 5532988154186

$ws.Range("A7").Value = 6
$ws.Range("C7").Value = 5532988220411

# --- Row heights: rows 2-3 grow from 60 to 75, new rows 4-7 also 75 ---
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(3).RowHeight = 75
$ws.Rows.Item(4).RowHeight = 75
$ws.Rows.Item(5).RowHeight = 75
$ws.Rows.Item(6).RowHeight = 75
$ws.Rows.Item(7).RowHeight = 75

# --- Update the active selection ---
$ws.Range("H5").Select()
